# Update column F ("dSF") values for rows 2-23 (excluding rows 13 and 19,
# which are unchanged) on the active worksheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    2  = -1
    3  = -1
    4  = -3
    5  = 3
    6  = -5
    7  = 2
    8  = 0
    9  = -1
    10 = -4
    11 = -3
    12 = 3
    14 = 2
    15 = -5
    16 = -1
    17 = -2
    18 = -1
    20 = -2
    21 = -4
    22 = -1
    23 = 6
}

foreach ($row in $updates.Keys) {
    $ws.Cells.Item($row, 6).Value = $updates[$row]
}
